$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 0.1764705882352941
$ws.Range("C2").Value = 0.6006191950464397
$ws.Range("J2").Value = 0.01238390092879257
$ws.Range("P2").Value = 0.1455108359133127
$ws.Range("S2").Value = 0.06501547987616099
# Row 3
$ws.Range("B3").Value = 0.01030927835051546
$ws.Range("C3").Value = 0.01030927835051546
$ws.Range("J3").Value = 0.02061855670103093
$ws.Range("P3").Value = 0.7628865979381443
$ws.Range("S3").Value = 0.1958762886597938
# Row 4
$ws.Range("P4").Value = 0.75
$ws.Range("S4").Value = 0.25
# Row 6
$ws.Range("B6").Value = 0.06132075471698113
$ws.Range("D6").Value = 0.01886792452830189
$ws.Range("F6").Value = 0.03773584905660377
$ws.Range("J6").Value = 0.2452830188679245
$ws.Range("O6").Value = 0.01886792452830189
$ws.Range("Q6").Value = 0.1415094339622641
$ws.Range("R6").Value = 0.09433962264150944
$ws.Range("S6").Value = 0.3820754716981132
# Row 7
$ws.Range("B7").Value = 0.1376146788990826
$ws.Range("D7").Value = 0.01834862385321101
$ws.Range("F7").Value = 0.04587155963302753
$ws.Range("J7").Value = 0.1100917431192661
$ws.Range("O7").Value = 0.02752293577981652
$ws.Range("Q7").Value = 0.1559633027522936
$ws.Range("R7").Value = 0.06880733944954129
$ws.Range("S7").Value = 0.4357798165137615
# Row 8
$ws.Range("B8").Value = 0.1425233644859813
$ws.Range("D8").Value = 0.01635514018691589
$ws.Range("F8").Value = 0.06542056074766354
$ws.Range("J8").Value = 0.1168224299065421
$ws.Range("O8").Value = 0.01401869158878505
$ws.Range("Q8").Value = 0.1401869158878505
$ws.Range("R8").Value = 0.09813084112149532
$ws.Range("S8").Value = 0.4065420560747663
# Row 9
$ws.Range("B9").Value = 0.1446808510638298
$ws.Range("D9").Value = 0.01702127659574468
$ws.Range("F9").Value = 0.03404255319148936
$ws.Range("J9").Value = 0.1021276595744681
$ws.Range("O9").Value = 0.01702127659574468
$ws.Range("Q9").Value = 0.1872340425531915
$ws.Range("R9").Value = 0.0851063829787234
$ws.Range("S9").Value = 0.4127659574468085
# Row 10
$ws.Range("B10").Value = 0.1041844577284372
$ws.Range("D10").Value = 0.02988898377455167
$ws.Range("F10").Value = 0.07941929974380871
$ws.Range("J10").Value = 0.09649871904355252
$ws.Range("O10").Value = 0.02049530315969257
$ws.Range("Q10").Value = 0.215200683176772
$ws.Range("R10").Value = 0.09479077711357814
$ws.Range("S10").Value = 0.3595217762596072
# Row 11
$ws.Range("G11").Value = 0.1402439024390244
$ws.Range("J11").Value = 0.08841463414634146
$ws.Range("K11").Value = 0.1829268292682927
$ws.Range("L11").Value = 0.5762195121951219
$ws.Range("S11").Value = 0.01219512195121951
# Row 12
$ws.Range("G12").Value = 0.7411167512690355
$ws.Range("J12").Value = 0.1776649746192893
$ws.Range("L12").Value = 0.04060913705583756
$ws.Range("S12").Value = 0.04060913705583756
# Row 13
$ws.Range("G13").Value = 0.75
$ws.Range("J13").Value = 0.1590909090909091
$ws.Range("S13").Value = 0.09090909090909091
# Row 14
$ws.Range("G14").Value = 0.75
$ws.Range("J14").Value = 0.25
# Row 15
$ws.Range("F15").Value = 0.01456310679611651
$ws.Range("H15").Value = 0.145631067961165
$ws.Range("I15").Value = 0.06796116504854369
$ws.Range("J15").Value = 0.3155339805825243
$ws.Range("K15").Value = 0.05825242718446602
$ws.Range("M15").Value = 0.02912621359223301
$ws.Range("O15").Value = 0.05825242718446602
$ws.Range("S15").Value = 0.3106796116504854
# Row 16
$ws.Range("F16").Value = 0.03539823008849557
$ws.Range("H16").Value = 0.1946902654867257
$ws.Range("I16").Value = 0.1194690265486726
$ws.Range("J16").Value = 0.3185840707964602
$ws.Range("K16").Value = 0.1327433628318584
$ws.Range("M16").Value = 0.02654867256637168
$ws.Range("O16").Value = 0.04867256637168142
$ws.Range("S16").Value = 0.1238938053097345
# Row 17
$ws.Range("F17").Value = 0.01435406698564593
$ws.Range("H17").Value = 0.1913875598086124
$ws.Range("I17").Value = 0.1004784688995215
$ws.Range("J17").Value = 0.3923444976076555
$ws.Range("K17").Value = 0.1028708133971292
$ws.Range("M17").Value = 0.02870813397129187
$ws.Range("N17").Value = 0.002392344497607655
$ws.Range("O17").Value = 0.05741626794258373
$ws.Range("S17").Value = 0.1100478468899522
# Row 18
$ws.Range("F18").Value = 0.00966183574879227
$ws.Range("H18").Value = 0.1594202898550725
$ws.Range("I18").Value = 0.106280193236715
$ws.Range("J18").Value = 0.4106280193236715
$ws.Range("K18").Value = 0.08695652173913043
$ws.Range("M18").Value = 0.00966183574879227
$ws.Range("N18").Value = 0.00966183574879227
$ws.Range("O18").Value = 0.08695652173913043
$ws.Range("S18").Value = 0.1207729468599034
# Row 19
$ws.Range("F19").Value = 0.01203852327447833
$ws.Range("H19").Value = 0.2014446227929374
$ws.Range("I19").Value = 0.1051364365971107
$ws.Range("J19").Value = 0.3627608346709471
$ws.Range("K19").Value = 0.1284109149277688
$ws.Range("M19").Value = 0.014446227929374
$ws.Range("N19").Value = 0.001605136436597111
$ws.Range("O19").Value = 0.06179775280898876
$ws.Range("S19").Value = 0.1123595505617977